# Generate Report for Handoff
# Updates the "Latest HO Xliff Generate Date" timestamps and sets the
# "Priority" column to "ht" for the report rows on the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$rows = @(7, 8, 9, 12, 13, 14)

foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-09-02 22:24:43"
    $wsZhCn.Range("H$r").Value = "2016-09-02 22:24:35"
    $wsDeDe.Range("H$r").Value = "2016-09-02 22:24:43"

    $wsZhCn.Range("E$r").Value = "ht"
    $wsDeDe.Range("E$r").Value = "ht"
}
